$wb = $excel.ActiveWorkbook

# Updated "想去人数" (want-to-go count) values for sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 570
$ws1.Range("F5").Value  = 285
$ws1.Range("F6").Value  = 1092
$ws1.Range("F7").Value  = 1433
$ws1.Range("F8").Value  = 587
$ws1.Range("F9").Value  = 108
$ws1.Range("F10").Value = 750
$ws1.Range("F11").Value = 70
$ws1.Range("F12").Value = 160
$ws1.Range("F15").Value = 1351
$ws1.Range("F17").Value = 102
$ws1.Range("F20").Value = 20
$ws1.Range("F21").Value = 651
$ws1.Range("F22").Value = 1003
$ws1.Range("F23").Value = 34
$ws1.Range("F24").Value = 220
$ws1.Range("F26").Value = 5840
$ws1.Range("F27").Value = 62
$ws1.Range("F31").Value = 14465
$ws1.Range("F32").Value = 1435
$ws1.Range("F33").Value = 211
$ws1.Range("F36").Value = 8549
$ws1.Range("F37").Value = 616
$ws1.Range("F38").Value = 4204
$ws1.Range("F40").Value = 359

# Updated "想去人数" (want-to-go count) values for sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 570
$ws4.Range("F5").Value  = 285
$ws4.Range("F6").Value  = 1092
$ws4.Range("F7").Value  = 1433
$ws4.Range("F8").Value  = 587
$ws4.Range("F9").Value  = 108
$ws4.Range("F10").Value = 750
$ws4.Range("F11").Value = 70
$ws4.Range("F12").Value = 160
$ws4.Range("F15").Value = 1351
$ws4.Range("F17").Value = 102
$ws4.Range("F21").Value = 20
$ws4.Range("F22").Value = 651
$ws4.Range("F24").Value = 1003
$ws4.Range("F25").Value = 34
$ws4.Range("F26").Value = 220
$ws4.Range("F29").Value = 5840
$ws4.Range("F30").Value = 62
$ws4.Range("F34").Value = 14465
$ws4.Range("F35").Value = 1435
$ws4.Range("F36").Value = 211
$ws4.Range("F39").Value = 8549
$ws4.Range("F40").Value = 616
$ws4.Range("F41").Value = 4204
$ws4.Range("F43").Value = 359

$wb.Save()
